{"js": "// Add detailed caIntegrator status to three bullet items under the\n// TRANSCEND section (caArray refresh, role-based permissions, single\n// sign-on), per the \"Added detailed caIntegrator status\" commit.\n\n// 1) \"Automatic caArray Refresh \u2013 Ongoing\"\n//    -> append the parenthetical CAINT ticket summary.\nconst autoRefreshResults = context.document.body.search(\n  \"Automatic caArray Refresh \\u2013 Ongoing\",\n  { matchCase: true }\n);\nautoRefreshResults.load(\"items\");\nawait context.sync();\n\nif (autoRefreshResults.items.length === 0) {\n  throw new Error('Could not find paragraph \"Automatic caArray Refresh \u2013 Ongoing\"');\n}\n\nautoRefreshResults.items[0].insertText(\n  \" (CAINT-1115, CAINT-1116, CAINT-1117, CAINT-1118 were completed.\" +\n    \" The remaining subtasks for CAINT-1094 should be completed by the end if Iteration I18.)\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n\n// 2) \"Role based Permissions \u2013 Ongoing\"\n//    -> append the iteration note.\nconst rolePermResults = context.document.body.search(\n  \"Role based Permissions \\u2013 Ongoing\",\n  { matchCase: true }\n);\nrolePermResults.load(\"items\");\nawait context.sync();\n\nif (rolePermResults.items.length === 0) {\n  throw new Error('Could not find paragraph \"Role based Permissions \u2013 Ongoing\"');\n}\n\nrolePermResults.items[0].insertText(\n  \" (CAINT-1107 being worked on in Iteration I17)\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n\n// 3) \"Single Sign-ON - Ongoing\"\n//    -> normalize the plain hyphen to an en dash (\"Single Sign-ON \u2013 Ongoing\").\nconst ssoResults = context.document.body.search(\"Single Sign-ON - Ongoing\", {\n  matchCase: true\n});\nssoResults.load(\"items\");\nawait context.sync();\n\nif (ssoResults.items.length === 0) {\n  throw new Error('Could not find paragraph \"Single Sign-ON - Ongoing\"');\n}\n\nssoResults.items[0].insertText(\n  \"Single Sign-ON \\u2013 Ongoing\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Add detailed caIntegrator status to three bullet items under TRANSCEND.\n\n$d = $word.ActiveDocument\n\n# 1) \"Automatic caArray Refresh \u2013 Ongoing\" -> append parenthetical ticket detail.\n$range = $d.Content\n$found = $range.Find.Execute(\"Automatic caArray Refresh \" + [char]0x2013 + \" Ongoing\")\nif (-not $found) {\n    throw \"Could not find paragraph 'Automatic caArray Refresh - Ongoing'\"\n}\n$range.Collapse(0)  # wdCollapseEnd\n$range.InsertAfter(\" (CAINT-1115, CAINT-1116, CAINT-1117, CAINT-1118 were completed.\" + `\n    \" The remaining subtasks for CAINT-1094 should be completed by the end if Iteration I18.)\")\n\n# 2) \"Role based Permissions \u2013 Ongoing\" -> append iteration note.\n$range = $d.Content\n$found = $range.Find.Execute(\"Role based Permissions \" + [char]0x2013 + \" Ongoing\")\nif (-not $found) {\n    throw \"Could not find paragraph 'Role based Permissions - Ongoing'\"\n}\n$range.Collapse(0)  # wdCollapseEnd\n$range.InsertAfter(\" (CAINT-1107 being worked on in Iteration I17)\")\n\n# 3) \"Single Sign-ON - Ongoing\" -> normalize the hyphen to an en dash.\n$range = $d.Content\n$found = $range.Find.Execute(\"Single Sign-ON - Ongoing\")\nif (-not $found) {\n    throw \"Could not find paragraph 'Single Sign-ON - Ongoing'\"\n}\n$range.Text = \"Single Sign-ON \" + [char]0x2013 + \" Ongoing\"\n"}
